# Apply updated transition-probability values for Tulsa_B sheet
# (team spec time commit pt2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2
$ws.Range("C2").Value = 0.518840579710145
$ws.Range("J2").Value = 0.01449275362318841
$ws.Range("O2").Value = 0.005797101449275362
$ws.Range("P2").Value = 0.1449275362318841
$ws.Range("S2").Value = 0.1159420289855072
# Row 3
$ws.Range("B3").Value = 0.0053475935828877
$ws.Range("C3").Value = 0.0374331550802139
$ws.Range("J3").Value = 0.0481283422459893
$ws.Range("P3").Value = 0.7272727272727273
$ws.Range("S3").Value = 0.1818181818181818
# Row 4
$ws.Range("J4").Value = 0.06818181818181818
$ws.Range("P4").Value = 0.6363636363636364
$ws.Range("S4").Value = 0.2954545454545455
# Row 6
$ws.Range("B6").Value = 0.04975124378109453
$ws.Range("D6").Value = 0.01492537313432836
$ws.Range("F6").Value = 0.03980099502487562
$ws.Range("J6").Value = 0.318407960199005
$ws.Range("O6").Value = 0.01990049751243781
$ws.Range("Q6").Value = 0.109452736318408
$ws.Range("R6").Value = 0.07462686567164178
$ws.Range("S6").Value = 0.373134328358209
# Row 7
$ws.Range("B7").Value = 0.1216730038022814
$ws.Range("D7").Value = 0.01520912547528517
$ws.Range("F7").Value = 0.03041825095057034
$ws.Range("J7").Value = 0.155893536121673
$ws.Range("O7").Value = 0.03041825095057034
$ws.Range("Q7").Value = 0.1863117870722434
$ws.Range("R7").Value = 0.06463878326996197
$ws.Range("S7").Value = 0.3954372623574144
# Row 8
$ws.Range("B8").Value = 0.09953703703703703
$ws.Range("D8").Value = 0.02314814814814815
$ws.Range("F8").Value = 0.07407407407407407
$ws.Range("J8").Value = 0.1111111111111111
$ws.Range("O8").Value = 0.03240740740740741
$ws.Range("Q8").Value = 0.1666666666666667
$ws.Range("R8").Value = 0.08796296296296297
$ws.Range("S8").Value = 0.4050925925925926
# Row 9
$ws.Range("B9").Value = 0.07981220657276995
$ws.Range("D9").Value = 0.02347417840375587
$ws.Range("F9").Value = 0.06103286384976526
$ws.Range("J9").Value = 0.1220657276995305
$ws.Range("O9").Value = 0.01408450704225352
$ws.Range("Q9").Value = 0.1924882629107981
$ws.Range("R9").Value = 0.08450704225352113
$ws.Range("S9").Value = 0.4225352112676056
# Row 10
$ws.Range("B10").Value = 0.1127348643006263
$ws.Range("D10").Value = 0.01530967292971468
$ws.Range("E10").Value = 0.0006958942240779402
$ws.Range("F10").Value = 0.05080027835768963
$ws.Range("J10").Value = 0.1217814892136395
$ws.Range("O10").Value = 0.01948503827418233
$ws.Range("Q10").Value = 0.2045929018789144
$ws.Range("R10").Value = 0.08907446068197634
$ws.Range("S10").Value = 0.3855254001391789
# Row 11
$ws.Range("G11").Value = 0.1648648648648649
$ws.Range("J11").Value = 0.05675675675675676
$ws.Range("K11").Value = 0.1972972972972973
$ws.Range("L11").Value = 0.5621621621621622
$ws.Range("S11").Value = 0.01891891891891892
# Row 12
$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.1745283018867924
$ws.Range("K12").Value = 0.01886792452830189
$ws.Range("L12").Value = 0.01886792452830189
$ws.Range("S12").Value = 0.03773584905660377
# Row 13
$ws.Range("G13").Value = 0.71875
$ws.Range("J13").Value = 0.21875
$ws.Range("S13").Value = 0.0625
# Row 15
$ws.Range("F15").Value = 0.02304147465437788
$ws.Range("H15").Value = 0.1658986175115207
$ws.Range("I15").Value = 0.06912442396313365
$ws.Range("J15").Value = 0.3179723502304148
$ws.Range("K15").Value = 0.05990783410138249
$ws.Range("M15").Value = 0.02304147465437788
$ws.Range("O15").Value = 0.04608294930875576
$ws.Range("S15").Value = 0.2949308755760369
# Row 16
$ws.Range("H16").Value = 0.1625615763546798
$ws.Range("I16").Value = 0.06403940886699508
$ws.Range("J16").Value = 0.477832512315271
$ws.Range("K16").Value = 0.1083743842364532
$ws.Range("M16").Value = 0.04926108374384237
$ws.Range("O16").Value = 0.04926108374384237
$ws.Range("S16").Value = 0.08866995073891626
# Row 17
$ws.Range("F17").Value = 0.02953586497890295
$ws.Range("H17").Value = 0.1877637130801688
$ws.Range("I17").Value = 0.08016877637130802
$ws.Range("J17").Value = 0.4367088607594937
$ws.Range("K17").Value = 0.109704641350211
$ws.Range("M17").Value = 0.01054852320675105
$ws.Range("O17").Value = 0.04641350210970464
$ws.Range("S17").Value = 0.09915611814345991
# Row 18
$ws.Range("F18").Value = 0.0184331797235023
$ws.Range("H18").Value = 0.1336405529953917
$ws.Range("I18").Value = 0.07834101382488479
$ws.Range("J18").Value = 0.4285714285714285
$ws.Range("K18").Value = 0.1612903225806452
$ws.Range("M18").Value = 0.02304147465437788
$ws.Range("N18").Value = 0.004608294930875576
$ws.Range("O18").Value = 0.05990783410138249
$ws.Range("S18").Value = 0.09216589861751152
# Row 19
$ws.Range("F19").Value = 0.01466275659824047
$ws.Range("H19").Value = 0.1825513196480938
$ws.Range("I19").Value = 0.09604105571847507
$ws.Range("J19").Value = 0.4032258064516129
$ws.Range("K19").Value = 0.1202346041055719
$ws.Range("M19").Value = 0.03079178885630499
$ws.Range("O19").Value = 0.05205278592375367
$ws.Range("S19").Value = 0.1004398826979472
